$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = 5
    4  = 17
    5  = 12
    6  = 15
    7  = 9
    8  = 8
    9  = 15
    10 = 6
    11 = 17
    12 = 18
    13 = 18
    14 = 13
    15 = 6
    17 = 13
    18 = 11
    19 = 14
    20 = 17
    21 = 7
    22 = 18
    23 = 11
    24 = 20
    25 = 5
    26 = 12
    27 = 10
    28 = 13
    29 = 8
    30 = 5
    32 = 6
    33 = 13
    34 = 9
    35 = 9
    36 = 11
    37 = 17
    39 = 15
    40 = 16
    41 = 12
    42 = 5
    43 = 14
    44 = 9
    45 = 20
    46 = 12
    47 = 20
    48 = 14
    49 = 8
    50 = 17
    51 = 13
    52 = 14
    53 = 19
    54 = 17
    56 = 19
    57 = 18
    58 = 5
    59 = 13
    60 = 16
    61 = 5
    62 = 11
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}
